# "Generate Report for Archive"
#
# The localization-status report is regenerated: the two files
# 66ec3868-ad0c-45fa-9a94-6068a8de1ec8 and 5261da22-23e2-4c5f-a60b-446bc987709a
# trade places (66ec3868 now sorts into row 4 / 5261da22 into row 5 on every
# sheet), their Status moves from "Ready for handoff" to "In Translation" for
# 255e04c5 and 66ec3868 (5261da22 becomes "Ready for handoff"), and their
# per-language handoff file / datetime columns are refreshed to point at the
# correct (own) artifact instead of the stale one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B3").Value = "In Translation"
$ov.Range("C3").Value = "In Translation"

$ov.Range("A4").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"
$ov.Range("D4").Value = "2016-33-18 20:33:10"

$ov.Range("A5").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
$ov.Range("D5").Value = "2016-32-18 20:32:01"

foreach ($hl in $ov.Hyperlinks) {
    if ($hl.Range.Row -eq 4 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
    }
    if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C3").Value = "In Translation"

$zh.Range("A4").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
$zh.Range("C4").Value = "In Translation"
$zh.Range("D4").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-18 20:33:07"

$zh.Range("A5").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
$zh.Range("D5").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.zh-cn.xlf"
$zh.Range("E5").Value = "2016-03-18 20:31:58"

foreach ($hl in $zh.Hyperlinks) {
    if ($hl.Range.Row -eq 4 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
    }
    if ($hl.Range.Row -eq 4 -and $hl.Range.Column -eq 4) {
        $hl.TextToDisplay = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.zh-cn.xlf"
    }
    if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
    }
    if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 4) {
        $hl.TextToDisplay = "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C3").Value = "In Translation"

$de.Range("A4").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
$de.Range("C4").Value = "In Translation"
$de.Range("D4").Value = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.de-de.xlf"
$de.Range("E4").Value = "2016-03-18 20:33:10"

$de.Range("A5").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
$de.Range("D5").Value = "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.de-de.xlf"
$de.Range("E5").Value = "2016-03-18 20:32:01"

foreach ($hl in $de.Hyperlinks) {
    if ($hl.Range.Row -eq 4 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.md"
    }
    if ($hl.Range.Row -eq 4 -and $hl.Range.Column -eq 4) {
        $hl.TextToDisplay = "66ec3868-ad0c-45fa-9a94-6068a8de1ec8.b926bccb71c10874a9e9887bbb707a64496b2e33.de-de.xlf"
    }
    if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "5261da22-23e2-4c5f-a60b-446bc987709a.md"
    }
    if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 4) {
        $hl.TextToDisplay = "5261da22-23e2-4c5f-a60b-446bc987709a.8cb54633748861568b40efef69ca42cfeb87fd10.de-de.xlf"
    }
}
